$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.442.30'
$ws.Range('E2').Value = '  -0.87%  '
$ws.Range('D3').Value = '3.897.55'
$ws.Range('E3').Value = '  +3.87%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '602.19'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '164.54'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.85%  '
$ws.Range('D7').Value = '3.896.46'
$ws.Range('E7').Value = '  +3.87%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('E9').Value = '  -1.92%  '
$ws.Range('E10').Value = '  -3.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.36'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.460'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.17%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '36.87'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.11%  '
$ws.Range('E14').Value = '  -1.45%  '
$ws.Range('D15').Value = '4.551.02'
$ws.Range('E15').Value = '  +3.90%  '
$ws.Range('D16').Value = '3.864.13'
$ws.Range('E16').Value = '  +3.35%  '
$ws.Range('D17').Value = '68.625.97'
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('E19').Value = '  -0.94%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.94'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.19'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.53%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '483.54'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.717'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000166'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +11.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.20'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.41%  '
$ws.Range('E26').Value = '  -1.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.03'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.09'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('E30').Value = '  -1.21%  '
$ws.Range('D31').Value = '4.050.69'
$ws.Range('E31').Value = '  +4.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.84'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.45%  '
$ws.Range('E33').Value = '  -2.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '31.86'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.70%  '
$ws.Range('D35').Value = '3.840.50'
$ws.Range('E35').Value = '  +4.04%  '
$ws.Range('E37').Value = '  +2.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.139'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.87'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.92%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('E41').Value = '  -1.81%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '437.17'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.02%  '
$ws.Range('E43').Value = '  -4.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '48.40'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.25%  '
$ws.Range('E45').Value = '  -0.48%  '
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.45'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.28%  '
$ws.Range('D48').Value = '2.833.83'
$ws.Range('E48').Value = '  +1.72%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '142.06'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '25.97'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +9.46%  '
$ws.Range('E51').Value = '  +0.15%  '
